$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "'09/21/2025"
$ws.Range("A20").ClearFormats()
$ws.Range("B20").Value = 0.1307131026099658
$ws.Range("C20").Value = 0.8692868973900342
